{"js": "// Justify (align \"both\" in OOXML terms) every body paragraph except the\n// title. The title paragraph keeps its existing centered alignment; all\n// other paragraphs get Word.Alignment.justified, matching <w:jc w:val=\"both\"/>\n// being inserted into each of their <w:pPr>.\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nfor (let i = 1; i < paragraphs.items.length; i++) {\n  paragraphs.items[i].alignment = Word.Alignment.justified;\n}\n\nawait context.sync();\n", "ps1": "# Justify (align \"both\" in OOXML terms) every body paragraph except the\n# title (paragraph 1, which keeps its existing centered alignment).\n# wdAlignParagraphJustify = 3 -> <w:jc w:val=\"both\"/>\n$d = $word.ActiveDocument\n\n$count = $d.Paragraphs.Count\nfor ($i = 2; $i -le $count; $i++) {\n    $d.Paragraphs.Item($i).Alignment = 3\n}\n"}
